$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1190.625
$ws.Range("I33").Value = 152.15384
$ws.Range("K33").Value = 152.15384
$ws.Range("M33").Value = 76.84616
$ws.Range("H41").Value = 1060.3846
$ws.Range("I41").Value = 985.25
$ws.Range("J41").Value = 1093.7778
$ws.Range("K41").Value = 985.25
$ws.Range("L41").Value = 1093.7778
$ws.Range("M41").Value = -545.25
$ws.Range("N41").Value = -1973.7778
$ws.Range("H76").Value = 4864.143
$ws.Range("I76").Value = 4762.25
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4762.25
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4447.25
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 4864.143
$ws.Range("I79").Value = 4762.25
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4762.25
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3670.25
$ws.Range("N79").Value = -7184
$ws.Range("H129").Value = 2159.1143
$ws.Range("J129").Value = 830.10205
$ws.Range("L129").Value = 2490.30615
$ws.Range("N129").Value = -12490.30615
$ws.Range("H132").Value = 4812273.5
$ws.Range("I132").Value = 5323896.5
$ws.Range("K132").Value = 15971689.5
$ws.Range("M132").Value = -15969159.5
$ws.Range("H135").Value = 1154.8125
$ws.Range("I135").Value = 322.86206
$ws.Range("J135").Value = 2424.6316
$ws.Range("K135").Value = 2905.75854
$ws.Range("L135").Value = 21821.6844
$ws.Range("M135").Value = -370.7585399999998
$ws.Range("N135").Value = -26891.6844
$ws.Range("H141").Value = 3015.4546
$ws.Range("I141").Value = 2855.7144
$ws.Range("J141").Value = 3295
$ws.Range("K141").Value = 8567.143199999999
$ws.Range("L141").Value = 9885
$ws.Range("M141").Value = -3387.143199999999
$ws.Range("N141").Value = -20245

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46317.047
$ws.Range("I2").Value = 948.0714
$ws.Range("K2").Value = 948.0714
$ws.Range("M2").Value = -835.0714
$ws.Range("H116").Value = 46317.047
$ws.Range("I116").Value = 948.0714
$ws.Range("K116").Value = 948.0714
$ws.Range("M116").Value = 1345.9286
$ws.Range("H132").Value = 2519.194
$ws.Range("I132").Value = 2545.6604
$ws.Range("J132").Value = 2419
$ws.Range("K132").Value = 7636.9812
$ws.Range("L132").Value = 7257
$ws.Range("M132").Value = -5106.9812
$ws.Range("N132").Value = -12317

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46317.047
$ws.Range("I3").Value = 948.0714
$ws.Range("K3").Value = 948.0714
$ws.Range("M3").Value = -834.0714
$ws.Range("H20").Value = 61476.824
$ws.Range("I20").Value = 69387.07000000001
$ws.Range("J20").Value = 2150
$ws.Range("K20").Value = 69387.07000000001
$ws.Range("L20").Value = 2150
$ws.Range("M20").Value = -69140.07000000001
$ws.Range("N20").Value = -2644
$ws.Range("H86").Value = 70530.75
$ws.Range("I86").Value = 86239.38
$ws.Range("J86").Value = 2460
$ws.Range("K86").Value = 86239.38
$ws.Range("L86").Value = 2460
$ws.Range("M86").Value = -85116.38
$ws.Range("N86").Value = -4706
$ws.Range("H89").Value = 70530.75
$ws.Range("I89").Value = 86239.38
$ws.Range("J89").Value = 2460
$ws.Range("K89").Value = 431196.9
$ws.Range("L89").Value = 12300
$ws.Range("M89").Value = -425580.9
$ws.Range("N89").Value = -23532
$ws.Range("H134").Value = 1822.1041
$ws.Range("I134").Value = 1612.4667
$ws.Range("J134").Value = 4966.6665
$ws.Range("K134").Value = 4837.4001
$ws.Range("L134").Value = 14899.9995
$ws.Range("M134").Value = -2302.4001
$ws.Range("N134").Value = -19969.9995

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5500
$ws.Range("J29").Value = 5500
$ws.Range("L29").Value = 5500
$ws.Range("N29").Value = -6086
$ws.Range("H31").Value = 2373.9216
$ws.Range("I31").Value = 605.2
$ws.Range("J31").Value = 4074.6155
$ws.Range("K31").Value = 605.2
$ws.Range("L31").Value = 4074.6155
$ws.Range("M31").Value = -310.2
$ws.Range("N31").Value = -4664.6155
$ws.Range("H34").Value = 2373.9216
$ws.Range("I34").Value = 605.2
$ws.Range("J34").Value = 4074.6155
$ws.Range("K34").Value = 605.2
$ws.Range("L34").Value = 4074.6155
$ws.Range("M34").Value = -403.2
$ws.Range("N34").Value = -4478.6155

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 8528.143
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 4500
$ws.Range("M22").Value = -4331
$ws.Range("H27").Value = 8528.143
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 4500
$ws.Range("M27").Value = -4398
$ws.Range("H34").Value = 692.0526
$ws.Range("J34").Value = 906.9286
$ws.Range("L34").Value = 2720.7858
$ws.Range("N34").Value = -2888.7858
$ws.Range("H58").Value = 2733.3333
$ws.Range("I58").Value = 2600
$ws.Range("K58").Value = 7800
$ws.Range("M58").Value = -7672
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 6000
$ws.Range("M125").Value = -1080
$ws.Range("H131").Value = 989.70966
$ws.Range("J131").Value = 1000.47253
$ws.Range("L131").Value = 3001.41759
$ws.Range("N131").Value = -13081.41759

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 124401
$ws.Range("I70").Value = 131613
$ws.Range("J70").Value = 9009
$ws.Range("K70").Value = 131613
$ws.Range("L70").Value = 9009
$ws.Range("M70").Value = -131343
$ws.Range("N70").Value = -9549
$ws.Range("H73").Value = 124401
$ws.Range("I73").Value = 131613
$ws.Range("J73").Value = 9009
$ws.Range("K73").Value = 131613
$ws.Range("L73").Value = 9009
$ws.Range("M73").Value = -130677
$ws.Range("N73").Value = -10881
$ws.Range("H102").Value = 3529.182
$ws.Range("I102").Value = 3380.1667
$ws.Range("K102").Value = 3380.1667
$ws.Range("M102").Value = -1758.1667

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2617.9
$ws.Range("J22").Value = 4999.6665
$ws.Range("L22").Value = 4999.6665
$ws.Range("N22").Value = -5589.6665
$ws.Range("H27").Value = 2617.9
$ws.Range("J27").Value = 4999.6665
$ws.Range("L27").Value = 4999.6665
$ws.Range("N27").Value = -5213.6665

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 48140.25
$ws.Range("J33").Value = 60914
$ws.Range("L33").Value = 60914
$ws.Range("N33").Value = -61414
$ws.Range("H36").Value = 48140.25
$ws.Range("J36").Value = 60914
$ws.Range("L36").Value = 60914
$ws.Range("N36").Value = -61414
$ws.Range("H37").Value = 13017.714
$ws.Range("J37").Value = 13799.667
$ws.Range("L37").Value = 13799.667
$ws.Range("N37").Value = -14205.667
$ws.Range("H56").Value = 3666.5
$ws.Range("I56").Value = 4000
$ws.Range("J56").Value = 3333
$ws.Range("K56").Value = 4000
$ws.Range("L56").Value = 3333
$ws.Range("M56").Value = -3286
$ws.Range("N56").Value = -4761
